$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55 (ALC)
$ws.Range("H55").Value = 435.7143
$ws.Range("I55").Value = 250.33333
$ws.Range("J55").Value = 574.75
$ws.Range("K55").Value = 250.33333
$ws.Range("L55").Value = 574.75
$ws.Range("M55").Value = -36.33332999999999
$ws.Range("N55").Value = -1002.75

# Row 74 (ALC)
$ws.Range("H74").Value = 86683
$ws.Range("I74").Value = 4019.6
$ws.Range("J74").Value = 500000
$ws.Range("K74").Value = 4019.6
$ws.Range("L74").Value = 500000
$ws.Range("M74").Value = -3083.6
$ws.Range("N74").Value = -501872

# Row 77 (ALC)
$ws.Range("H77").Value = 86683
$ws.Range("I77").Value = 4019.6
$ws.Range("J77").Value = 500000
$ws.Range("K77").Value = 20098
$ws.Range("L77").Value = 2500000
$ws.Range("M77").Value = -15418
$ws.Range("N77").Value = -2509360

# Row 86 (ALC)
$ws.Range("H86").Value = 4934.375
$ws.Range("I86").Value = 4639.2856
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 4639.2856
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -3516.2856
$ws.Range("N86").Value = -9246

# Row 89 (ALC)
$ws.Range("H89").Value = 4934.375
$ws.Range("I89").Value = 4639.2856
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 23196.428
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -17580.428
$ws.Range("N89").Value = -46232

# Row 100 (ALC)
$ws.Range("H100").Value = 3863
$ws.Range("I100").Value = 697.25
$ws.Range("J100").Value = 6395.6
$ws.Range("K100").Value = 697.25
$ws.Range("L100").Value = 6395.6
$ws.Range("M100").Value = -156.25
$ws.Range("N100").Value = -7477.6

# Row 103 (ALC)
$ws.Range("H103").Value = 766
$ws.Range("I103").Value = 498.5
$ws.Range("J103").Value = 899.75
$ws.Range("K103").Value = 1495.5
$ws.Range("L103").Value = 2699.25
$ws.Range("M103").Value = -909.5
$ws.Range("N103").Value = -3871.25

# Row 104 (ALC)
$ws.Range("H104").Value = 259.8
$ws.Range("I104").Value = 259.8
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 779.4000000000001
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 967.5999999999999

# Row 111 (ALC)
$ws.Range("H111").Value = 748.8
$ws.Range("I111").Value = 773.5
$ws.Range("J111").Value = 650
$ws.Range("K111").Value = 2320.5
$ws.Range("L111").Value = 1950
$ws.Range("M111").Value = 746.5
$ws.Range("N111").Value = -8084

# Row 112 (ALC)
$ws.Range("H112").Value = 1607.7391
$ws.Range("I112").Value = 1209.8
$ws.Range("J112").Value = 1718.2778
$ws.Range("K112").Value = 3629.4
$ws.Range("L112").Value = 5154.8334
$ws.Range("M112").Value = -2521.4
$ws.Range("N112").Value = -7370.8334

# Row 132 (ALC)
$ws.Range("H132").Value = 2241
$ws.Range("I132").Value = 2249.7693
$ws.Range("J132").Value = 2203
$ws.Range("K132").Value = 6749.3079
$ws.Range("L132").Value = 6609
$ws.Range("M132").Value = -4219.3079
$ws.Range("N132").Value = -11669

# Row 137 (ALC)
$ws.Range("H137").Value = 2699.6924
$ws.Range("I137").Value = 1877.6666
$ws.Range("J137").Value = 4549.25
$ws.Range("K137").Value = 5632.9998
$ws.Range("L137").Value = 13647.75
$ws.Range("M137").Value = -3082.9998
$ws.Range("N137").Value = -18747.75

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 2662.889
$ws.Range("I61").Value = 2662.889
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2662.889
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2450.889

# Row 122 (ARM)
$ws.Range("H122").Value = 2300
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -9700

# Row 132 (ARM)
$ws.Range("H132").Value = 1464.4546
$ws.Range("I132").Value = 1464.4546
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4393.3638
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1863.3638

# Row 136 (ARM)
$ws.Range("H136").Value = 2662.889
$ws.Range("I136").Value = 2662.889
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7988.667
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5438.667

$ws = $wb.Worksheets.Item("BSM")
# Row 41 (BSM)
$ws.Range("H41").Value = 199999
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 199999
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 199999
$ws.Range("N41").Value = -200775

# Row 134 (BSM)
$ws.Range("H134").Value = 8032.357
$ws.Range("I134").Value = 9639.315000000001
$ws.Range("J134").Value = 4639.8887
$ws.Range("K134").Value = 28917.945
$ws.Range("L134").Value = 13919.6661
$ws.Range("M134").Value = -26382.945
$ws.Range("N134").Value = -18989.6661

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (CRP)
$ws.Range("H7").Value = 160.07692
$ws.Range("I7").Value = 95.25
$ws.Range("J7").Value = 263.8
$ws.Range("K7").Value = 95.25
$ws.Range("L7").Value = 263.8
$ws.Range("M7").Value = 17.75
$ws.Range("N7").Value = -489.8

# Row 62 (CRP)
$ws.Range("H62").Value = 4600.4
$ws.Range("I62").Value = 1002.5
$ws.Range("J62").Value = 6999
$ws.Range("K62").Value = 1002.5
$ws.Range("L62").Value = 6999
$ws.Range("M62").Value = -378.5
$ws.Range("N62").Value = -8247

# Row 65 (CRP)
$ws.Range("H65").Value = 4600.4
$ws.Range("I65").Value = 1002.5
$ws.Range("J65").Value = 6999
$ws.Range("K65").Value = 5012.5
$ws.Range("L65").Value = 34995
$ws.Range("M65").Value = -1892.5
$ws.Range("N65").Value = -41235

# Row 103 (CRP)
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()

# Row 122 (CRP)
$ws.Range("H122").Value = 1400
$ws.Range("I122").Value = 1400
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4200
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1750

# Row 134 (CRP)
$ws.Range("H134").Value = 5224.875
$ws.Range("I134").Value = 4780
$ws.Range("J134").Value = 5966.3335
$ws.Range("K134").Value = 14340
$ws.Range("L134").Value = 17899.0005
$ws.Range("M134").Value = -11805
$ws.Range("N134").Value = -22969.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (CUL)
$ws.Range("H12").Value = 94.333336
$ws.Range("I12").Value = 65
$ws.Range("J12").Value = 100.2
$ws.Range("K12").Value = 195
$ws.Range("L12").Value = 300.6
$ws.Range("M12").Value = -22
$ws.Range("N12").Value = -646.6

# Row 80 (CUL)
$ws.Range("H80").Value = 20000
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 60000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -59064

# Row 83 (CUL)
$ws.Range("H83").Value = 20000
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 180000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -175320

# Row 92 (CUL)
$ws.Range("H92").Value = 874.75
$ws.Range("I92").Value = 874.75
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2624.25
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1376.25
$ws.Range("N92").ClearContents()

# Row 112 (CUL)
$ws.Range("H112").Value = 500
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -392

$ws = $wb.Worksheets.Item("GSM")
# Row 53 (GSM)
$ws.Range("H53").Value = 29143
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 29143
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 29143
$ws.Range("N53").Value = -30405

# Row 58 (GSM)
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# Row 80 (GSM)
$ws.Range("H80").Value = 3031.5557
$ws.Range("I80").Value = 2399
$ws.Range("J80").Value = 3110.625
$ws.Range("K80").Value = 2399
$ws.Range("L80").Value = 3110.625
$ws.Range("M80").Value = -1401
$ws.Range("N80").Value = -5106.625

# Row 83 (GSM)
$ws.Range("H83").Value = 3031.5557
$ws.Range("I83").Value = 2399
$ws.Range("J83").Value = 3110.625
$ws.Range("K83").Value = 11995
$ws.Range("L83").Value = 15553.125
$ws.Range("M83").Value = -7003
$ws.Range("N83").Value = -25537.125

# Row 132 (GSM)
$ws.Range("H132").Value = 2180.8333
$ws.Range("I132").Value = 2017
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 6051
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -3521
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Range("H61").Value = 8943.6
$ws.Range("I61").Value = 8943.6
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 8943.6
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -8741.6

# Row 87 (LTW)
$ws.Range("H87").Value = 39997
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 39997
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 39997
$ws.Range("N87").Value = -42243

# Row 90 (LTW)
$ws.Range("H90").Value = 39997
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 39997
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 119991
$ws.Range("N90").Value = -131223

# Row 113 (LTW)
$ws.Range("H113").Value = 8943.6
$ws.Range("I113").Value = 8943.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 8943.6
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -6773.6

# Row 130 (LTW)
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 132 (LTW)
$ws.Range("H132").Value = 22500.25
$ws.Range("I132").Value = 17500.75
$ws.Range("J132").Value = 27499.75
$ws.Range("K132").Value = 52502.25
$ws.Range("L132").Value = 82499.25
$ws.Range("M132").Value = -49972.25
$ws.Range("N132").Value = -87559.25

# Row 136 (LTW)
$ws.Range("H136").Value = 2487.3635
$ws.Range("I136").Value = 2487.3635
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7462.0905
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4912.0905
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 4 (WVR)
$ws.Range("H4").Value = 17225
$ws.Range("I4").Value = 30450
$ws.Range("J4").Value = 4000
$ws.Range("K4").Value = 30450
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = -30337
$ws.Range("N4").Value = -4226

# Row 132 (WVR)
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -470

# Row 135 (WVR)
$ws.Range("H135").Value = 31398.8
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 31398.8
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 31398.8
$ws.Range("N135").Value = -41538.8
